# Apply the Tue Sep 12 06:25:36 UTC 2023 cryptos-list refresh.
# Coin / Link / Price / Volume(1h) columns are B / C / D / E; row 1 is the header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.843.78"
$ws.Range("E2").Value = "  -0.39%  "
# Row 3
$ws.Range("D3").Value = "1.586.20"
# Row 4
$ws.Range("E4").Value = "  +0.07%  "
# Row 5
$ws.Range("D5").Value = "'209.99"
$ws.Range("E5").Value = "  -1.27%  "
# Row 6
$ws.Range("E6").Value = "  +0.05%  "
# Row 7
$ws.Range("D7").Value = "'0.479"
$ws.Range("E7").Value = "  -3.68%  "
# Row 8
$ws.Range("E8").Value = "  -0.97%  "
# Row 9
$ws.Range("D9").Value = "'0.0616"
$ws.Range("E9").Value = "  -0.35%  "
# Row 10
$ws.Range("D10").Value = "'18.02"
$ws.Range("E10").Value = "  -2.15%  "
# Row 11
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  -0.34%  "
# Row 12
$ws.Range("D12").Value = "1.806.57"
$ws.Range("E12").Value = "  -2.15%  "
# Row 13
$ws.Range("D13").Value = "1.584.21"
$ws.Range("E13").Value = "  -2.30%  "
# Row 14
$ws.Range("E14").Value = "  -2.76%  "
# Row 15
$ws.Range("D15").Value = "'0.509"
$ws.Range("E15").Value = "  -3.00%  "
# Row 16
$ws.Range("D16").Value = "25.832.56"
$ws.Range("E16").Value = "  -0.51%  "
# Row 17
$ws.Range("D17").Value = "0.0₃0722"
$ws.Range("E17").Value = "  -2.13%  "
# Row 18
$ws.Range("D18").Value = "'59.83"
$ws.Range("E18").Value = "  -3.04%  "
# Row 19
$ws.Range("E19").Value = "  +0.03%  "
# Row 20
$ws.Range("D20").Value = "'191.53"
$ws.Range("E20").Value = "  -0.56%  "
# Row 21
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = "  -1.71%  "
# Row 22
$ws.Range("D22").Value = "'9.36"
$ws.Range("E22").Value = "  -1.82%  "
# Row 23
$ws.Range("E23").Value = "  -1.76%  "
# Row 24
$ws.Range("D24").Value = "'0.132"
$ws.Range("E24").Value = "  -0.99%  "
# Row 25
$ws.Range("D25").Value = "'141.98"
$ws.Range("E25").Value = "  -1.26%  "
# Row 26
$ws.Range("E26").Value = "  +0.03%  "
# Row 27
$ws.Range("D27").Value = "'1.70"
$ws.Range("E27").Value = "  -1.09%  "
# Row 28
$ws.Range("D28").Value = "'15.09"
$ws.Range("E28").Value = "  -1.05%  "
# Row 29
$ws.Range("E29").Value = "  -3.01%  "
# Row 30
$ws.Range("E30").Value = "  -5.71%  "
# Row 31
$ws.Range("D31").Value = "'0.0470"
$ws.Range("E31").Value = "  -1.79%  "
# Row 32
$ws.Range("E32").Value = "  -0.36%  "
# Row 33
$ws.Range("E33").Value = "  -2.33%  "
# Row 34
$ws.Range("E34").Value = "  +0.05%  "
# Row 35
$ws.Range("E35").Value = "  -2.41%  "
# Row 36
$ws.Range("D36").Value = "1.101.50"
$ws.Range("E36").Value = "  -2.20%  "
# Row 37
$ws.Range("E37").Value = "  -0.01%  "
# Row 38
$ws.Range("D38").Value = "'2.35"
$ws.Range("E38").Value = "  -2.20%  "
# Row 39
$ws.Range("D39").Value = "'0.503"
$ws.Range("E39").Value = "  -2.39%  "
# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0150"
$ws.Range("E40").Value = "  -2.11%  "
# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.776"
$ws.Range("E41").Value = "  -8.10%  "
# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.817"
$ws.Range("E42").Value = "  +7.88%  "
# Row 43
$ws.Range("E43").Value = "  +2.13%  "
# Row 44
$ws.Range("D44").Value = "'93.90"
$ws.Range("E44").Value = "  -3.94%  "
# Row 45
$ws.Range("D45").Value = "1.720.47"
$ws.Range("E45").Value = "  -2.10%  "
# Row 46
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  -0.96%  "
# Row 47
$ws.Range("E47").Value = "  -1.07%  "
# Row 48
$ws.Range("D48").Value = "'53.20"
$ws.Range("E48").Value = "  -1.70%  "
# Row 49
$ws.Range("D49").Value = "'0.0508"
$ws.Range("E49").Value = "  -1.77%  "
# Row 50
$ws.Range("E50").Value = "  -0.93%  "
# Row 51
$ws.Range("E51").Value = "  -0.07%  "
